$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("Z74").Value = "작성테스트"
$ws.Range("Z74").Font.Name = "돋움"
$ws.Range("Z74").Font.Size = 10
$ws.Range("Z74").Font.Color = 0
Write-Output "done"
